$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 407, pushing the existing rows 407..464 down to 408..465.
$ws.Rows("407:407").Insert()

# Populate the newly inserted row 407 with the new price-report record.
$ws.Cells.Item(407, 1).Value = 10
$ws.Cells.Item(407, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(407, 3).Value = "La Araucanía"
$ws.Cells.Item(407, 4).Value = 44816
$ws.Cells.Item(407, 5).Value = 9
$ws.Cells.Item(407, 6).Value = 100112024
$ws.Cells.Item(407, 7).Value = "Choclo"
$ws.Cells.Item(407, 8).Value = "Dulce o Americano"
$ws.Cells.Item(407, 9).Value = "Primera"
$ws.Cells.Item(407, 10).Value = 65
$ws.Cells.Item(407, 11).Value = 24000
$ws.Cells.Item(407, 12).Value = 25000
$ws.Cells.Item(407, 13).Value = 24462
$ws.Cells.Item(407, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(407, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(407, 16).Value = 349
$ws.Cells.Item(407, 17).Value = 70
$ws.Cells.Item(407, 18).Value = "Hortaliza"
